$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (new string order observed first: 04104017374, 23/04/2021, 12:30)
$ws.Range("B4").Value = "ssurgwsoadev4-oci.opc.oracleoutsourcing.com"
$ws.Range("C4").Value = "https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/cc/ClaimCenter.do"
$ws.Range("C4").Style = "Normal"
$ws.Range("F4").Value = "'04104017374"
$ws.Range("H4").Value = "'23/04/2021"
$ws.Range("I4").Value = "'12:30"

# Row 2 (new strings: 04104017412, 26/04/2021)
$ws.Range("F2").Value = "'04104017412"
$ws.Range("H2").Value = "'26/04/2021"

# Row 3 (reuses same new strings as row 2)
$ws.Range("F3").Value = "'04104017412"
$ws.Range("H3").Value = "'26/04/2021"

# Update selection to match target view state
$ws.Range("H4").Select()
